$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 46770
$ws.Cells.Item(2, 5).Value = 5110
$ws.Cells.Item(2, 6).Value = 5110
$ws.Cells.Item(2, 7).Value = 4802
$ws.Cells.Item(2, 8).Value = 3546
$ws.Cells.Item(2, 9).Value = 3494
$ws.Cells.Item(2, 10).Value = 51
$ws.Cells.Item(2, 11).Value = 38283
$ws.Cells.Item(2, 12).Value = 21198
$ws.Cells.Item(2, 13).Value = 17085
$ws.Cells.Item(2, 14).Value = 16374
$ws.Cells.Item(2, 15).Value = 711
$ws.Cells.Item(2, 16).Value = 886
$ws.Cells.Item(2, 17).Value = 4658
$ws.Cells.Item(2, 18).Value = -2776
$ws.Cells.Item(2, 19).Value = 77
$ws.Cells.Item(2, 20).Value = 1872
$ws.Cells.Item(2, 21).Value = 2786
$ws.Cells.Item(2, 22).Value = 12086
$ws.Cells.Item(2, 23).Value = 10.93
$ws.Cells.Item(2, 24).Value = 7.58
$ws.Cells.Item(2, 25).Value = 22.94
$ws.Cells.Item(2, 26).Value = 9.76
$ws.Cells.Item(2, 27).Value = 124.08
$ws.Cells.Item(2, 28).Value = 2016.45
$ws.Cells.Item(2, 29).Value = 19722
$ws.Cells.Item(2, 30).Value = 31.59
$ws.Cells.Item(2, 31).Value = 97722
$ws.Cells.Item(2, 32).Value = 6.38
$ws.Cells.Item(2, 33).Value = 4000
$ws.Cells.Item(2, 34).Value = 0.64
$ws.Cells.Item(2, 35).Value = 19.21
$ws.Cells.Item(2, 36).Value = 15618197

# Row 3
$ws.Cells.Item(3, 4).Value = 53285
$ws.Cells.Item(3, 5).Value = 6841
$ws.Cells.Item(3, 6).Value = 6841
$ws.Cells.Item(3, 7).Value = 6448
$ws.Cells.Item(3, 8).Value = 4704
$ws.Cells.Item(3, 9).Value = 4604
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 42146
$ws.Cells.Item(3, 12).Value = 20993
$ws.Cells.Item(3, 13).Value = 21153
$ws.Cells.Item(3, 14).Value = 20373
$ws.Cells.Item(3, 15).Value = 780
$ws.Cells.Item(3, 16).Value = 886
$ws.Cells.Item(3, 17).Value = 6378
$ws.Cells.Item(3, 18).Value = -3470
$ws.Cells.Item(3, 19).Value = -2315
$ws.Cells.Item(3, 20).Value = 3017
$ws.Cells.Item(3, 21).Value = 3361
$ws.Cells.Item(3, 22).Value = 10612
$ws.Cells.Item(3, 23).Value = 12.84
$ws.Cells.Item(3, 24).Value = 8.83
$ws.Cells.Item(3, 25).Value = 25.05
$ws.Cells.Item(3, 26).Value = 11.7
$ws.Cells.Item(3, 27).Value = 99.24
$ws.Cells.Item(3, 28).Value = 2442.24
$ws.Cells.Item(3, 29).Value = 25982
$ws.Cells.Item(3, 30).Value = 40.41
$ws.Cells.Item(3, 31).Value = 121586
$ws.Cells.Item(3, 32).Value = 8.64
$ws.Cells.Item(3, 33).Value = 5500
$ws.Cells.Item(3, 34).Value = 0.52
$ws.Cells.Item(3, 35).Value = 20.04
$ws.Cells.Item(3, 36).Value = 15618197

# Row 4
$ws.Cells.Item(4, 4).Value = 60941
$ws.Cells.Item(4, 5).Value = 8809
$ws.Cells.Item(4, 6).Value = 8809
$ws.Cells.Item(4, 7).Value = 7527
$ws.Cells.Item(4, 8).Value = 5792
$ws.Cells.Item(4, 9).Value = 5682
$ws.Cells.Item(4, 10).Value = 110
$ws.Cells.Item(4, 11).Value = 45022
$ws.Cells.Item(4, 12).Value = 18813
$ws.Cells.Item(4, 13).Value = 26208
$ws.Cells.Item(4, 14).Value = 25324
$ws.Cells.Item(4, 15).Value = 885
$ws.Cells.Item(4, 16).Value = 886
$ws.Cells.Item(4, 17).Value = 7134
$ws.Cells.Item(4, 18).Value = -4055
$ws.Cells.Item(4, 19).Value = -3670
$ws.Cells.Item(4, 20).Value = 3315
$ws.Cells.Item(4, 21).Value = 3819
$ws.Cells.Item(4, 22).Value = 7943
$ws.Cells.Item(4, 23).Value = 14.46
$ws.Cells.Item(4, 24).Value = 9.5
$ws.Cells.Item(4, 25).Value = 24.87
$ws.Cells.Item(4, 26).Value = 13.29
$ws.Cells.Item(4, 27).Value = 71.78
$ws.Cells.Item(4, 28).Value = 2976.43
$ws.Cells.Item(4, 29).Value = 32070
$ws.Cells.Item(4, 30).Value = 26.72
$ws.Cells.Item(4, 31).Value = 151131
$ws.Cells.Item(4, 32).Value = 5.67
$ws.Cells.Item(4, 33).Value = 7500
$ws.Cells.Item(4, 34).Value = 0.88
$ws.Cells.Item(4, 35).Value = 22.13
$ws.Cells.Item(4, 36).Value = 15618197

# Row 5
$ws.Cells.Item(5, 4).Value = 61051
$ws.Cells.Item(5, 5).Value = 9300
$ws.Cells.Item(5, 6).Value = 9300
$ws.Cells.Item(5, 7).Value = 8611
$ws.Cells.Item(5, 8).Value = 6183
$ws.Cells.Item(5, 9).Value = 6064
$ws.Cells.Item(5, 10).Value = 119
$ws.Cells.Item(5, 11).Value = 47785
$ws.Cells.Item(5, 12).Value = 16967
$ws.Cells.Item(5, 13).Value = 30818
$ws.Cells.Item(5, 14).Value = 29962
$ws.Cells.Item(5, 15).Value = 855
$ws.Cells.Item(5, 16).Value = 886
$ws.Cells.Item(5, 17).Value = 7355
$ws.Cells.Item(5, 18).Value = -3338
$ws.Cells.Item(5, 19).Value = -3511
$ws.Cells.Item(5, 20).Value = 2842
$ws.Cells.Item(5, 21).Value = 4513
$ws.Cells.Item(5, 22).Value = 6041
$ws.Cells.Item(5, 23).Value = 15.23
$ws.Cells.Item(5, 24).Value = 10.13
$ws.Cells.Item(5, 25).Value = 21.94
$ws.Cells.Item(5, 26).Value = 13.33
$ws.Cells.Item(5, 27).Value = 55.06
$ws.Cells.Item(5, 28).Value = 3546.1
$ws.Cells.Item(5, 29).Value = 34226
$ws.Cells.Item(5, 30).Value = 34.74
$ws.Cells.Item(5, 31).Value = 178815
$ws.Cells.Item(5, 32).Value = 6.65
$ws.Cells.Item(5, 33).Value = 9000
$ws.Cells.Item(5, 34).Value = 0.76
$ws.Cells.Item(5, 35).Value = 24.89
$ws.Cells.Item(5, 36).Value = 15618197

# Row 6
$ws.Cells.Item(6, 4).Value = 67475
$ws.Cells.Item(6, 5).Value = 10392
$ws.Cells.Item(6, 6).Value = 10392
$ws.Cells.Item(6, 7).Value = 9560
$ws.Cells.Item(6, 8).Value = 6923
$ws.Cells.Item(6, 9).Value = 6827
$ws.Cells.Item(6, 10).Value = 95
$ws.Cells.Item(6, 11).Value = 52759
$ws.Cells.Item(6, 12).Value = 16819
$ws.Cells.Item(6, 13).Value = 35940
$ws.Cells.Item(6, 14).Value = 35121
$ws.Cells.Item(6, 15).Value = 820
$ws.Cells.Item(6, 16).Value = 886
$ws.Cells.Item(6, 17).Value = 8171
$ws.Cells.Item(6, 18).Value = -4303
$ws.Cells.Item(6, 19).Value = -3780
$ws.Cells.Item(6, 20).Value = 3563
$ws.Cells.Item(6, 21).Value = 4608
$ws.Cells.Item(6, 22).Value = 4814
$ws.Cells.Item(6, 23).Value = 15.4
$ws.Cells.Item(6, 24).Value = 10.26
$ws.Cells.Item(6, 25).Value = 20.98
$ws.Cells.Item(6, 26).Value = 13.77
$ws.Cells.Item(6, 27).Value = 46.8
$ws.Cells.Item(6, 28).Value = 4132.2
$ws.Cells.Item(6, 29).Value = 38534
$ws.Cells.Item(6, 30).Value = 28.57
$ws.Cells.Item(6, 31).Value = 209600
$ws.Cells.Item(6, 32).Value = 5.25
$ws.Cells.Item(6, 33).Value = 9250
$ws.Cells.Item(6, 34).Value = 0.84
$ws.Cells.Item(6, 35).Value = 22.72
$ws.Cells.Item(6, 36).Value = 15618197

# Row 7
$ws.Cells.Item(7, 4).Value = 76141
$ws.Cells.Item(7, 5).Value = 11757
$ws.Cells.Item(7, 7).Value = 11112
$ws.Cells.Item(7, 8).Value = 8050
$ws.Cells.Item(7, 9).Value = 7954
$ws.Cells.Item(7, 11).Value = 60955
$ws.Cells.Item(7, 12).Value = 18484
$ws.Cells.Item(7, 13).Value = 42471
$ws.Cells.Item(7, 14).Value = 41529
$ws.Cells.Item(7, 16).Value = 888
$ws.Cells.Item(7, 17).Value = 9443
$ws.Cells.Item(7, 18).Value = -4736
$ws.Cells.Item(7, 19).Value = -2514
$ws.Cells.Item(7, 20).Value = 3523
$ws.Cells.Item(7, 21).Value = 5683
$ws.Cells.Item(7, 23).Value = 15.44
$ws.Cells.Item(7, 24).Value = 10.57
$ws.Cells.Item(7, 25).Value = 20.75
$ws.Cells.Item(7, 26).Value = 14.16
$ws.Cells.Item(7, 27).Value = 43.52
$ws.Cells.Item(7, 29).Value = 44891
$ws.Cells.Item(7, 30).Value = 30.05
$ws.Cells.Item(7, 31).Value = 247844
$ws.Cells.Item(7, 32).Value = 5.44
$ws.Cells.Item(7, 33).Value = 9757
$ws.Cells.Item(7, 34).Value = 0.72
$ws.Cells.Item(7, 35).Value = 19.16

# Row 8
$ws.Cells.Item(8, 4).Value = 84484
$ws.Cells.Item(8, 5).Value = 12868
$ws.Cells.Item(8, 7).Value = 12141
$ws.Cells.Item(8, 8).Value = 8840
$ws.Cells.Item(8, 9).Value = 8749
$ws.Cells.Item(8, 11).Value = 68914
$ws.Cells.Item(8, 12).Value = 19477
$ws.Cells.Item(8, 13).Value = 49437
$ws.Cells.Item(8, 14).Value = 48421
$ws.Cells.Item(8, 16).Value = 889
$ws.Cells.Item(8, 17).Value = 10354
$ws.Cells.Item(8, 18).Value = -4232
$ws.Cells.Item(8, 19).Value = -1990
$ws.Cells.Item(8, 20).Value = 3936
$ws.Cells.Item(8, 21).Value = 6216
$ws.Cells.Item(8, 23).Value = 15.23
$ws.Cells.Item(8, 24).Value = 10.46
$ws.Cells.Item(8, 25).Value = 19.45
$ws.Cells.Item(8, 26).Value = 13.57
$ws.Cells.Item(8, 27).Value = 39.4
$ws.Cells.Item(8, 29).Value = 49377
$ws.Cells.Item(8, 30).Value = 25.48
$ws.Cells.Item(8, 31).Value = 288977
$ws.Cells.Item(8, 32).Value = 4.35
$ws.Cells.Item(8, 33).Value = 10555
$ws.Cells.Item(8, 34).Value = 0.84
$ws.Cells.Item(8, 35).Value = 18.84

# Row 9
$ws.Cells.Item(9, 4).Value = 91366
$ws.Cells.Item(9, 5).Value = 14342
$ws.Cells.Item(9, 7).Value = 13659
$ws.Cells.Item(9, 8).Value = 9923
$ws.Cells.Item(9, 9).Value = 9827
$ws.Cells.Item(9, 11).Value = 77368
$ws.Cells.Item(9, 12).Value = 19871
$ws.Cells.Item(9, 13).Value = 57498
$ws.Cells.Item(9, 14).Value = 56362
$ws.Cells.Item(9, 16).Value = 889
$ws.Cells.Item(9, 17).Value = 11490
$ws.Cells.Item(9, 18).Value = -4094
$ws.Cells.Item(9, 19).Value = -1978
$ws.Cells.Item(9, 20).Value = 3923
$ws.Cells.Item(9, 21).Value = 7403
$ws.Cells.Item(9, 23).Value = 15.7
$ws.Cells.Item(9, 24).Value = 10.86
$ws.Cells.Item(9, 25).Value = 18.76
$ws.Cells.Item(9, 26).Value = 13.57
$ws.Cells.Item(9, 27).Value = 34.56
$ws.Cells.Item(9, 29).Value = 55465
$ws.Cells.Item(9, 30).Value = 22.68
$ws.Cells.Item(9, 31).Value = 336366
$ws.Cells.Item(9, 32).Value = 3.74
$ws.Cells.Item(9, 33).Value = 11180
$ws.Cells.Item(9, 34).Value = 0.89
$ws.Cells.Item(9, 35).Value = 17.77
